# [2023-08-11] Major Updates on ROLODEX
# Rebuilds the Sheet1 "Rolodex" header row with many new CRM-style
# columns, drops the old two-row sample/hyperlink data, restyles the
# sheet (wrap-text + date columns), and updates sheet defaults.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) New header cells H1:AC1 - reuse the existing "header" look (fill +
#    border) already applied to G1 by copying its format across, then
#    fill in the new column titles.
# ---------------------------------------------------------------------
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:AC1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$headers = @(
    "H1","LinkedIn URL",
    "I1","Facebook URL",
    "J1","Twitter URL",
    "K1","Instagram URL",
    "L1","Industry",
    "M1","NAICS Code",
    "N1","Employees",
    "O1","Annual Revenue",
    "P1","Type",
    "Q1","OPT OUT",
    "R1","Billing Street",
    "S1","Billing City",
    "T1","Billing State",
    "U1","Billing Postal Code",
    "V1","Billing Country",
    "W1","Shipping Street",
    "X1","Shipping City",
    "Y1","Shipping State",
    "Z1","Shipping Postal Code",
    "AA1","Shipping Country",
    "AB1","Description",
    "AC1","Created Date & Time"
)
for ($i = 0; $i -lt $headers.Count; $i += 2) {
    $ws.Range($headers[$i]).Value = $headers[$i + 1]
}

# ---------------------------------------------------------------------
# 2) Drop the old two sample rows (their hyperlinks + "Hyperlink" cell
#    style go with them).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Range("A2:G3").Clear()
$wb.Styles.Item("Hyperlink").Delete()

# ---------------------------------------------------------------------
# 3) New blank body rows 2-6: row 2 gets a slightly custom height, R2
#    gets wrap-text, and the "Created Date & Time" column (AC) gets a
#    date/time number format down through row 6.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 13.8

$ws.Range("R2").WrapText = $true
$ws.Range("AC2").NumberFormat = "m/d/yy h:mm"

$ws.Range("AC2").Copy() | Out-Null
$ws.Range("AC3:AC6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Column widths - best-fit equivalents for the new layout.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.72
$ws.Columns.Item(2).ColumnWidth = 25.5
$ws.Columns.Item(3).ColumnWidth = 14.72
$ws.Columns.Item(4).ColumnWidth = 12.61
$ws.Columns.Item(5).ColumnWidth = 13.28
$ws.Columns.Item(6).ColumnWidth = 13.72
$ws.Columns.Item(8).ColumnWidth = 10.72
$ws.Columns.Item(9).ColumnWidth = 12.06
$ws.Columns.Item(10).ColumnWidth = 9.83
$ws.Columns.Item(11).ColumnWidth = 12.06
$ws.Columns.Item(12).ColumnWidth = 6.72
$ws.Columns.Item(13).ColumnWidth = 9.83
$ws.Columns.Item(14).ColumnWidth = 8.94
$ws.Columns.Item(15).ColumnWidth = 13.28
$ws.Columns.Item(17).ColumnWidth = 7.61
$ws.Columns.Item(18).ColumnWidth = 10.28
$ws.Columns.Item(19).ColumnWidth = 8.5
$ws.Columns.Item(20).ColumnWidth = 9.61
$ws.Columns.Item(21).ColumnWidth = 15.17
$ws.Columns.Item(22).ColumnWidth = 11.83
$ws.Columns.Item(23).ColumnWidth = 12.28
$ws.Columns.Item(24).ColumnWidth = 10.39
$ws.Columns.Item(25).ColumnWidth = 11.5
$ws.Columns.Item(26).ColumnWidth = 17.17
$ws.Columns.Item(27).ColumnWidth = 13.83
$ws.Columns.Item(28).ColumnWidth = 9.39
$ws.Columns.Item(29).ColumnWidth = 17.39

# ---------------------------------------------------------------------
# 5) Selection + used range.
# ---------------------------------------------------------------------
$ws.Range("A2:AC6").Select()
